# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback
# XLIFF files were generated/processed, reflecting a newer report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 23:20:23"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 23:20:18"
$wsZhCn.Range("K2").Value = "2016-09-05 23:20:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 23:20:23"
$wsDeDe.Range("K2").Value = "2016-09-05 23:20:53"
